$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 111618089
$ws.Range("B2").Value = 96348
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "30"
$ws.Range("J2").Value = "plantor/tuvor"
$ws.Range("K2").Value = "blomning"
$ws.Range("L2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("P2").Value = "A 32649, Heda, Sm"
$ws.Range("Q2").Value = 580617.6201989455
$ws.Range("R2").Value = 6415136.627037819
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = "Kalmar"
$ws.Range("U2").Value = "Västervik"
$ws.Range("V2").Value = "Småland"
$ws.Range("W2").Value = "Gamleby"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2023-08-21"
$ws.Range("Z2").Value = "00:00"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2023-08-21"
$ws.Range("AB2").Value = "00:00"
$ws.Range("AC2").Value = "1 blomma"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AF2").Value = ""
$ws.Range("AG2").Value = $false
$ws.Range("AT2").Value = ""
$ws.Range("AW2").Value = "Magnus Kasselstrand"
$ws.Range("AX2").Value = "Magnus Kasselstrand"
$ws.Range("AY2").Value = ""
$ws.Range("A3").Value = 111618070
$ws.Range("B3").Value = 96348
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "15"
$ws.Range("J3").Value = "plantor/tuvor"
$ws.Range("K3").Value = "blomning"
$ws.Range("L3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("P3").Value = "A 32649, Heda, Sm"
$ws.Range("Q3").Value = 580592.470229132
$ws.Range("R3").Value = 6415141.442167919
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Kalmar"
$ws.Range("U3").Value = "Västervik"
$ws.Range("V3").Value = "Småland"
$ws.Range("W3").Value = "Gamleby"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-08-21"
$ws.Range("Z3").Value = "00:00"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-08-21"
$ws.Range("AB3").Value = "00:00"
$ws.Range("AC3").Value = "1 blomma"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AF3").Value = ""
$ws.Range("AG3").Value = $false
$ws.Range("AT3").Value = ""
$ws.Range("AW3").Value = "Magnus Kasselstrand"
$ws.Range("AX3").Value = "Magnus Kasselstrand"
$ws.Range("AY3").Value = ""
$ws.Range("A4").Value = 111618144
$ws.Range("B4").Value = 96348
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "2"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("P4").Value = "A 32649, Heda, Sm"
$ws.Range("Q4").Value = 580620.6996611424
$ws.Range("R4").Value = 6415142.541277731
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Kalmar"
$ws.Range("U4").Value = "Västervik"
$ws.Range("V4").Value = "Småland"
$ws.Range("W4").Value = "Gamleby"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-08-21"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-08-21"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AC4").Value = ""
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AF4").Value = ""
$ws.Range("AG4").Value = $false
$ws.Range("AT4").Value = ""
$ws.Range("AW4").Value = "Magnus Kasselstrand"
$ws.Range("AX4").Value = "Magnus Kasselstrand"
$ws.Range("AY4").Value = ""
$ws.Range("A5").Value = 111618039
$ws.Range("B5").Value = 93388
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 2180
$ws.Range("F5").Value = "Blåmossa"
$ws.Range("G5").Value = "Leucobryum glaucum"
$ws.Range("H5").Value = "(Hedw.) Ångstr."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("P5").Value = "A 32649, Heda, Sm"
$ws.Range("Q5").Value = 580599.6803078586
$ws.Range("R5").Value = 6415233.627682217
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Kalmar"
$ws.Range("U5").Value = "Västervik"
$ws.Range("V5").Value = "Småland"
$ws.Range("W5").Value = "Gamleby"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-08-21"
$ws.Range("Z5").Value = "00:00"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-08-21"
$ws.Range("AB5").Value = "00:00"
$ws.Range("AC5").Value = ""
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AF5").Value = ""
$ws.Range("AG5").Value = $false
$ws.Range("AT5").Value = ""
$ws.Range("AW5").Value = "Magnus Kasselstrand"
$ws.Range("AX5").Value = "Magnus Kasselstrand"
$ws.Range("AY5").Value = ""
$ws.Range("A6").Value = 111618078
$ws.Range("B6").Value = 96348
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 220787
$ws.Range("F6").Value = "Knärot"
$ws.Range("G6").Value = "Goodyera repens"
$ws.Range("H6").Value = "(L.) R. Br."
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = "plantor/tuvor"
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("P6").Value = "A 32649, Heda, Sm"
$ws.Range("Q6").Value = 580612.1009209087
$ws.Range("R6").Value = 6415119.491031807
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Kalmar"
$ws.Range("U6").Value = "Västervik"
$ws.Range("V6").Value = "Småland"
$ws.Range("W6").Value = "Gamleby"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-08-21"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-08-21"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AC6").Value = ""
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AF6").Value = ""
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = "Magnus Kasselstrand"
$ws.Range("AX6").Value = "Magnus Kasselstrand"
$ws.Range("AY6").Value = ""
$ws.Range("A7").Value = 111618046
$ws.Range("B7").Value = 93388
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 2180
$ws.Range("F7").Value = "Blåmossa"
$ws.Range("G7").Value = "Leucobryum glaucum"
$ws.Range("H7").Value = "(Hedw.) Ångstr."
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("P7").Value = "A 32649, Heda, Sm"
$ws.Range("Q7").Value = 580591.6383206119
$ws.Range("R7").Value = 6415156.322361182
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Kalmar"
$ws.Range("U7").Value = "Västervik"
$ws.Range("V7").Value = "Småland"
$ws.Range("W7").Value = "Gamleby"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2023-08-21"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2023-08-21"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AC7").Value = ""
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AF7").Value = ""
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = "Magnus Kasselstrand"
$ws.Range("AX7").Value = "Magnus Kasselstrand"
$ws.Range("AY7").Value = ""
$ws.Range("A9").Value = 111618056
$ws.Range("B9").Value = 96348
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "15"
$ws.Range("J9").Value = "plantor/tuvor"
$ws.Range("K9").Value = "blomning"
$ws.Range("L9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("P9").Value = "A 32649, Heda, Sm"
$ws.Range("Q9").Value = 580582.6881743574
$ws.Range("R9").Value = 6415124.22061418
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Kalmar"
$ws.Range("U9").Value = "Västervik"
$ws.Range("V9").Value = "Småland"
$ws.Range("W9").Value = "Gamleby"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-08-21"
$ws.Range("Z9").Value = "00:00"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2023-08-21"
$ws.Range("AB9").Value = "00:00"
$ws.Range("AC9").Value = "2 blommor"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AF9").Value = ""
$ws.Range("AG9").Value = $false
$ws.Range("AT9").Value = ""
$ws.Range("AW9").Value = "Magnus Kasselstrand"
$ws.Range("AX9").Value = "Magnus Kasselstrand"
$ws.Range("AY9").Value = ""